$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the old hyperlinks (A8 single cell + A2:A7 range) so the
#    relationship-id numbering starts clean for the new hyperlinks.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2. Shift the existing column-A content (and its per-row style) one
#    column to the right, into column B. Column A becomes blank but
#    keeps the style that was on each row (header style on row 1,
#    "looks-like-a-link" blue/underline style on the others).
# ------------------------------------------------------------------
$ws.Range("A1:A8").Cut($ws.Range("B1:B8"))
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3. Rewrite column B's text content for every row.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "li"

$ws.Range("B2").Value = "https://drive.google.com/open?id=13UNUhXirv7JB072YEreOERu_AuWlsAi2"
$ws.Range("B3").Value = "https://drive.google.com/file/d/1Y7IjsvIsigeld0jvKxJpAHEzq7rP4lSE/view"
$ws.Range("B4").Value = "https://drive.google.com/file/d/1Y7IjsvIsigeld0jvKxJpAHEzq7rP4lSE/view"
$ws.Range("B5").Value = "https://drive.google.com/file/d/1Y7IjsvIsigeld0jvKxJpAHEzq7rP4lSE/view"
$ws.Range("B6").Value = "https://drive.google.com/open?id=13UNUhXirv7JB072YEreOERu_AuWlsAi6"
$ws.Range("B7").Value = "https://drive.google.com/open?id=13UNUhXirv7JB072YEreOERu_AuWlsAi7"
$ws.Range("B8").Value = "https://drive.google.com/open?id=13UNUhXirv7JB072YEreOERu_AuWlsAi8"

# New 9th row.
$ws.Range("B9").Value = "https://drive.google.com/open?id=13UNUhXirv7JB072YEreOERu_AuWlsAi9"

# Give the new row's A cell the same "blank but styled" look as the
# rest of column A by copying the format from an already-styled cell.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 4. B6:B9 look like links (same underline / hyperlink-blue font) but
#    are not real hyperlinks. Get that formatting by adding a
#    throw-away hyperlink and immediately deleting it again - this
#    reuses Excel's built-in "Hyperlink" style without leaving a
#    relationship behind, and it keeps it off the rId numbering used
#    by the real hyperlinks added afterwards.
# ------------------------------------------------------------------
foreach ($addr in "B6", "B7", "B8", "B9") {
    $ws.Hyperlinks.Add($ws.Range($addr), $ws.Range($addr).Value())
    $ws.Hyperlinks.Item($ws.Hyperlinks.Count).Delete()
}

# ------------------------------------------------------------------
# 5. Real hyperlinks, added last so they land on rId1/rId2/rId3.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "https://drive.google.com/open?id=13UNUhXirv7JB072YEreOERu_AuWlsAi2")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://drive.google.com/file/d/1Y7IjsvIsigeld0jvKxJpAHEzq7rP4lSE/view")
$ws.Hyperlinks.Add($ws.Range("B4:B5"), "https://drive.google.com/file/d/1Y7IjsvIsigeld0jvKxJpAHEzq7rP4lSE/view", "", "", "https://drive.google.com/file/d/1Y7IjsvIsigeld0jvKxJpAHEzq7rP4lSE/view")

# ------------------------------------------------------------------
# 6. Column sizing / view state.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 67.5
$excel.ActiveWindow.Zoom = 80
$ws.Range("B3:B5").Select()
